# ajustado restaurar backup com todos os campos da base
#
# Reconstructs the "restore backup" edit:
#  - conta_corrente: "Conta Padrão" account renamed/replaced by "ITAU"
#  - cartao_credito: selection cursor moved
#  - four new reference sheets added: estabelecimentos, contas, categorias, cartoes

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. conta_corrente ("Conta Padrão" -> "ITAU" on the conta column)
# ---------------------------------------------------------------------------
$contaCorrente = $wb.Worksheets.Item("conta_corrente")
$contaCorrente.Range("F2").Value = "ITAU"
$contaCorrente.Range("F3").Value = "ITAU"
$contaCorrente.Range("F4").Value = "ITAU"
$contaCorrente.Range("D3").Select()

# ---------------------------------------------------------------------------
# 2. cartao_credito (no data change, cursor moved)
# ---------------------------------------------------------------------------
$cartaoCredito = $wb.Worksheets.Item("cartao_credito")
$cartaoCredito.Range("F2").Select()

# ---------------------------------------------------------------------------
# Helper: copy the header look (bold font + border + centered alignment)
# from the existing conta_corrente header row onto a destination range.
# ---------------------------------------------------------------------------
function Copy-HeaderFormat($destRange) {
    $contaCorrente.Range("A1").Copy() | Out-Null
    $destRange.PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------------
# 3. estabelecimentos
# ---------------------------------------------------------------------------
$estabelecimentos = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$estabelecimentos.Name = "estabelecimentos"

$estabelecimentos.Range("A1").Value = "nome"
$estabelecimentos.Range("B1").Value = "descricao"
Copy-HeaderFormat($estabelecimentos.Range("A1:B1"))

$estabelecimentos.Range("A2").Value = "MEU EMPREGO"
$estabelecimentos.Range("A3").Value = "PAGUE MENOS"

$estabelecimentos.Range("A2:A3").Select()

# ---------------------------------------------------------------------------
# 4. contas
# ---------------------------------------------------------------------------
$contas = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$contas.Name = "contas"

$contas.Range("A1").Value = "nome"
$contas.Range("B1").Value = "descricao"
$contas.Range("C1").Value = "se_banco"
$contas.Range("D1").Value = "se_banco_nome"
$contas.Range("E1").Value = "se_banco_agencia"
$contas.Range("F1").Value = "se_banco_conta"
Copy-HeaderFormat($contas.Range("A1:F1"))

$contas.Range("A2").Value = "ITAU"
$contas.Range("C2").Value = $true
$contas.Range("D2").Value = "ITAU"

$contas.Range("E2").Select()

# ---------------------------------------------------------------------------
# 5. categorias
# ---------------------------------------------------------------------------
$categorias = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$categorias.Name = "categorias"

$categorias.Range("A1").Value = "nome"
$categorias.Range("B1").Value = "descricao"
$categorias.Range("C1").Value = "meta"
Copy-HeaderFormat($categorias.Range("A1:C1"))

$categorias.Range("A2").Value = "SALÁRIO"
$categorias.Range("C2").Value = 1000

$categorias.Range("A3").Value = "FARMÁCIA"
$categorias.Range("C3").Value = 0

$categorias.Range("A4").Value = "SORVETE"

$categorias.Range("C4").Select()

# ---------------------------------------------------------------------------
# 6. cartoes
# ---------------------------------------------------------------------------
$cartoes = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$cartoes.Name = "cartoes"

$cartoes.Range("A1").Value = "nome"
$cartoes.Range("B1").Value = "descricao"
$cartoes.Range("C1").Value = "bandeira"
$cartoes.Range("D1").Value = "ultimos_4_digitos"
$cartoes.Range("E1").Value = "dia_vencimento"
Copy-HeaderFormat($cartoes.Range("A1:E1"))

$cartoes.Range("A2").Value = "CARTAO PRINCIPAL"
$cartoes.Range("C2").Value = "OUTRO"
$cartoes.Range("E2").Value = 16

$cartoes.Range("E3").Select()

# ---------------------------------------------------------------------------
# Leave the user back on the first sheet, matching the saved view.
# ---------------------------------------------------------------------------
$contaCorrente.Activate()
$contaCorrente.Range("D3").Select()
